$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update totwAppearances for Julian Ryerson (row 18): 2 -> 3
$ws.Range("I18").Value = 3

# Remove the row for "Luca Reggiani" (row 21); this shifts the
# "Gregor Kobel" row (previously row 22) up to become row 21,
# shrinking the used range from A1:DL22 to A1:DL21.
$ws.Rows.Item(21).Delete()

# Update totwAppearances for Gregor Kobel (now row 21): 2 -> 3
$ws.Range("I21").Value = 3
